$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serials 46057-46059 => 2026-02-04..06)
$rows = @(
    @{ Row = 8;  Date = 46057; Station = "四方坪站"; C = 10059.94;  D = 8962.02;  E = 3780.64;              F = 445 },
    @{ Row = 9;  Date = 46057; Station = "高岭站";   C = 4439.49;   D = 4058.07;  E = 1192.83;              F = 165 },
    @{ Row = 10; Date = 46058; Station = "四方坪站"; C = 9987.98;   D = 9013.81;  E = 3792.07;              F = 433 },
    @{ Row = 11; Date = 46058; Station = "高岭站";   C = 3781.86;   D = 3383.53;  E = 1021.33;              F = 144 },
    @{ Row = 12; Date = 46059; Station = "四方坪站"; C = 11868.77;  D = 10982.1;  E = 4340.9399999999996;   F = 494 },
    @{ Row = 13; Date = 46059; Station = "高岭站";   C = 5170.7700000000004; D = 4577.26; E = 1385.29;      F = 170 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Station
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F

    # Copy styles from the row above (row-2 block, e.g. rows 6/7) to keep formatting consistent
    $srcRow = $row - 2
    $ws.Cells.Item($srcRow, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Cells.Item($srcRow, 3).Copy() | Out-Null
    $ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($srcRow, 4).Copy() | Out-Null
    $ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($srcRow, 5).Copy() | Out-Null
    $ws.Cells.Item($row, 5).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($srcRow, 6).Copy() | Out-Null
    $ws.Cells.Item($row, 6).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
